# [EDIT] trelo = trello aangepast in presentatie
# Fix the misspelled "Trelo" -> "Trello" in the "Planning/Trelo" bullet.

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -like "*Trelo*") {
                [void]$tr.Replace("Trelo", "Trello")
            }
        }
    }
}
